$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (from G1) onto the new H1 header cell, then set its value.
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$ws.Cells.Item(1, 8).Value = "Save"

# Fill the new "Save" column (H2:H7) with 0 for each data row.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
